$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = " "

$ws.Range("F13").Value = "wind"
$ws.Range("G13").Value = "Attack *1  Variation 20%"
$ws.Range("H13").Value = 30
$ws.Range("I13").Value = "on row of enemy"

$ws.Range("B15").Value = "Defence Boost"

$ws.Range("F17").Value = "Quick Sand"
$ws.Range("G17").Value = "Attack *.75 Variation 20%"
$ws.Range("I17").Value = "Front collum of enemy"
$ws.Range("J17").Value = "Decrease Speed By 20% 2 turns"

$ws.Range("F18").Value = "Earthquake"
$ws.Range("G18").Value = "Attack *1.25 Variation 20%"
$ws.Range("H18").Value = 60
$ws.Range("I18").Value = "All Enemy"

$ws.Range("G24").Select() | Out-Null
